# Word COM-interop script implementing the "wordsmithing" edits described
# in the commit. Several paragraphs of green (18A303) review-response text
# get their wording revised; a couple of paragraph-mark colors are
# normalised (removed/realigned) and one "Done" note changes from a
# slightly-off green (2CEE0E) to the canonical green (18A303).

$d = $word.ActiveDocument

$GREEN = 238360     # wdColor value for RGB 18A303 (BGR packed)
$AUTOC = -16777216  # wdColorAutomatic

function Set-ParaMarkAutoKeepGreenRuns($para) {
    # Clears the paragraph-mark (pilcrow) color override, while keeping the
    # visible run text green.
    $pStart = $para.Range.Start
    $pEnd = $para.Range.End
    $para.Range.Font.Color = $AUTOC
    if ($pEnd - 1 -gt $pStart) {
        $textRange = $d.Range($pStart, $pEnd - 1)
        $textRange.Font.Color = $GREEN
    }
}

function Replace-Text($para, [string]$oldText, [string]$newText) {
    $full = $para.Range.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        Write-Host "WARNING: text not found: $oldText"
        return
    }
    $start = $para.Range.Start + $idx
    $end = $start + $oldText.Length
    $rng = $d.Range($start, $end)
    $rng.Text = $newText
}

# ---------------------------------------------------------------------
# Paragraph: "Thank you for requesting that we incorporate our response
# from the first round into the manuscript. We included our reasoning to
# section 4.2"
# ---------------------------------------------------------------------
$p27 = $d.Paragraphs.Item(28)
Replace-Text $p27 "included our reasoning to section 4.2" "added our L shell dependence reasoning to section 4.2."
Set-ParaMarkAutoKeepGreenRuns $p27

# ---------------------------------------------------------------------
# Paragraph: "We added a paragraph explaining the imperfect normalization
# to section 4.2"  (paragraph-mark color 2CEE0E -> 18A303)
# ---------------------------------------------------------------------
$p33 = $d.Paragraphs.Item(34)
Replace-Text $p33 "a paragraph explaining the imperfect normalization to section 4.2" "a paragraph explaining the imperfect normalization to section 4.2."
$p33.Format.Font.Color = $GREEN

# ---------------------------------------------------------------------
# Paragraph: "Thank you for pointing this out, we made another attempt to
# clarify this section. ..."
# ---------------------------------------------------------------------
$p38 = $d.Paragraphs.Item(39)
$oldP38 = "Thank you for pointing this out, we made another attempt to clarify this section. Now we point out that the CDF distributions are similar since they both rapidly decrease towards a gradually-decreasing shoulder starting at 20 km in LEO which roughly scales to around 200 km at the magnetic equator."
$newP38 = "Thank you for pointing this out and we made another attempt to clarify this sentence. Besides mentioning that they are qualitatively similar, we now point out that both curves F(s) flatten out; the LEO curve flattens out around 20 km and the equatorial curve around 200 km separation. These separations correspond well with the magnetic field scaling from LEO to the magnetic equator."
Replace-Text $p38 $oldP38 $newP38
Set-ParaMarkAutoKeepGreenRuns $p38

# Wrap the bookmark text as in the original edit.
$full38 = $p38.Range.Text
$bookmarkText = "correspond well with the magnetic field scaling from LEO to the magnetic equator"
$bIdx = $full38.IndexOf($bookmarkText)
if ($bIdx -ge 0) {
    $bStart = $p38.Range.Start + $bIdx
    $bEnd = $bStart + $bookmarkText.Length
    $d.Bookmarks.Add("__DdeLink__47_2073078054", $d.Range($bStart, $bEnd))
}

# ---------------------------------------------------------------------
# Paragraph: "... Done" (run color 2CEE0E -> 18A303)
# ---------------------------------------------------------------------
$p42 = $d.Paragraphs.Item(43)
$rng42 = $p42.Range
$rng42.Find.ClearFormatting()
$rng42.Find.Replacement.ClearFormatting()
$rng42.Find.Replacement.Font.Color = $GREEN
$rng42.Find.Execute("Done", $true, $false, $false, $false, $false, $true, 1, $false, "Done", 2)

# ---------------------------------------------------------------------
# Paragraph: "When we addressed the 3rd moderate clarification we decided
# to point out the shoulder in both CDF curves so we believe that leaving
# it as CDF is more appropriate in this case."
# ---------------------------------------------------------------------
$p46 = $d.Paragraphs.Item(47)
$oldP46 = "point out the shoulder in both CDF curves so we believe that leaving it as CDF is more appropriate in this case. "
$newP46 = "point out the separation at which both CDF curves flatten out. So we believe that leaving it as CDF is more appropriate in this case. "
Replace-Text $p46 $oldP46 $newP46

Write-Host "Edits applied."
